# Add a new worksheet "ODI Batting Extra" as the last sheet in the workbook,
# matching the data scraped for extra batting fields.

$wb = $excel.ActiveWorkbook

# Add the sheet after the last existing sheet so it becomes sheet 4 (end of tab list).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Reuse the existing header style (bold, centered, thin border) from the
# "ODI Bowling" sheet so we don't introduce a brand-new style entry.
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBowling.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# Header row
$ws.Cells.Item(1,1).Value = "MATCH_CODE"
$ws.Cells.Item(1,2).Value = "BATTING_POSITION"
$ws.Cells.Item(1,3).Value = "NUM_4"
$ws.Cells.Item(1,4).Value = "NUM_6"
$ws.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Cells.Item(1,6).Value = "MAN_OF_MATCH"

# Column A (MATCH_CODE) holds text-formatted match codes for every data row.
$ws.Range("A2:A8").NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "4108"
$ws.Cells.Item(3,1).Value = "4115"
$ws.Cells.Item(4,1).Value = "4123"
$ws.Cells.Item(5,1).Value = "4125"
$ws.Cells.Item(6,1).Value = "4166"
$ws.Cells.Item(7,1).Value = "4167"
$ws.Cells.Item(8,1).Value = "4168"

# Column B (BATTING_POSITION) is a real number where known, blank otherwise.
$ws.Cells.Item(3,2).Value = 10
$ws.Cells.Item(4,2).Value = 9
$ws.Cells.Item(5,2).Value = 9
$ws.Cells.Item(6,2).Value = 9
$ws.Cells.Item(8,2).Value = 10

# Columns C (NUM_4) and D (NUM_6) are text-formatted digit counts where known.
$ws.Range("C3:D6").NumberFormat = "@"
$ws.Range("C8:D8").NumberFormat = "@"

$ws.Cells.Item(3,3).Value = "0"
$ws.Cells.Item(3,4).Value = "1"

$ws.Cells.Item(4,3).Value = "0"
$ws.Cells.Item(4,4).Value = "0"

$ws.Cells.Item(5,3).Value = "1"
$ws.Cells.Item(5,4).Value = "0"

$ws.Cells.Item(6,3).Value = "2"
$ws.Cells.Item(6,4).Value = "1"

$ws.Cells.Item(8,3).Value = "0"
$ws.Cells.Item(8,4).Value = "0"

# Column E (PERCENT_RUNS_OF_TOTAL) is a text percentage string where known.
$ws.Range("E3:E6").NumberFormat = "@"
$ws.Range("E8:E8").NumberFormat = "@"

$ws.Cells.Item(3,5).Value = "2.96%"
$ws.Cells.Item(4,5).Value = "1.52%"
$ws.Cells.Item(5,5).Value = "3.24%"
$ws.Cells.Item(6,5).Value = "8.88%"
$ws.Cells.Item(8,5).Value = "2.09%"

# Column F (MAN_OF_MATCH) is plain text "NO" for every data row.
$ws.Cells.Item(2,6).Value = "NO"
$ws.Cells.Item(3,6).Value = "NO"
$ws.Cells.Item(4,6).Value = "NO"
$ws.Cells.Item(5,6).Value = "NO"
$ws.Cells.Item(6,6).Value = "NO"
$ws.Cells.Item(7,6).Value = "NO"
$ws.Cells.Item(8,6).Value = "NO"

# Rows 2 and 7 leave BATTING_POSITION/NUM_4/NUM_6/PERCENT_RUNS_OF_TOTAL blank
# (matching the scraped source which had no data for those matches), so
# mark those cells as text-typed blanks.
$ws.Range("B2:E2").NumberFormat = "@"
$ws.Range("B2:E2").Value = ""
$ws.Range("B7:E7").NumberFormat = "@"
$ws.Range("B7:E7").Value = ""

$ws.Range("A1").Select() | Out-Null
